# "end of 4th week" - fill in progress % for week's tasks (D2:D11) and
# fill in the per-student grade/legend column (B14:B22), then leave the
# selection where the author ended up (B27), scrolled down to row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: "Степень выполнения" (percent complete) for tasks 1-10 ---
$ws.Range("D2").Value  = 0.5
$ws.Range("D3").Value  = 1
$ws.Range("D4").Value  = 0.9
$ws.Range("D5").Value  = 1
$ws.Range("D6").Value  = 0.8
$ws.Range("D7").Value  = 0.7
$ws.Range("D8").Value  = 0.7
$ws.Range("D9").Value  = 0.5
$ws.Range("D10").Value = 0.3
$ws.Range("D11").Value = 0

# --- Column B (rows 14-22): per-student grade / legend table ---
# Note: entry order matters so the new shared strings land at the same
# indices as the authored file ("5-"=31, "4-"=32, "позже"=33).
$ws.Range("B14").Value = "5-"
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 5
$ws.Range("B18").Value = "4-"
$ws.Range("B17").Value = "позже"
$ws.Range("B19").Value = "позже"
$ws.Range("B20").Value = "5-"
$ws.Range("B21").Value = 4
$ws.Range("B22").Value = 4

# --- View state: leave selection on B27, scrolled so row 10 is on top ---
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B27").Select()
